$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Login sheet - cosmetic resize only (row height / column width touch-up)
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A1:B2").EntireRow.RowHeight = 15
$wsLogin.Range("A1").EntireColumn.ColumnWidth = 21.71

# ---------------------------------------------------------------------------
# UserGroup sheet - cosmetic resize only
# ---------------------------------------------------------------------------
$wsUserGroup = $wb.Worksheets.Item("UserGroup")
$wsUserGroup.Range("A1:A2").EntireRow.RowHeight = 15
$wsUserGroup.Range("A1").EntireColumn.ColumnWidth = 17.86

# ---------------------------------------------------------------------------
# KnowledgeBase sheet - update description text, resize, reselect
# ---------------------------------------------------------------------------
$wsKB = $wb.Worksheets.Item("KnowledgeBase")
$wsKB.Range("C2").Value = "test"
$wsKB.Range("A1:C2").EntireRow.RowHeight = 15
$wsKB.Range("A1").EntireColumn.ColumnWidth = 26.59
$wsKB.Range("B1").EntireColumn.ColumnWidth = 32.57
$wsKB.Range("C1").EntireColumn.ColumnWidth = 35.58
$wsKB.Range("B12").Select()

# ---------------------------------------------------------------------------
# Alerts sheet - D2 becomes a real number, G2 keeps its value but is
# re-formatted as General, resize, tab / selection moves to RecordCreditNote
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("Alerts")
$wsAlerts.Range("D2").NumberFormat = "General"
$wsAlerts.Range("D2").Value = 10
$wsAlerts.Range("G2").NumberFormat = "General"

$wsAlerts.Range("A1:G2").EntireRow.RowHeight = 15
$wsAlerts.Range("A1").EntireColumn.ColumnWidth = 29.86
$wsAlerts.Range("B1").EntireColumn.ColumnWidth = 15.15
$wsAlerts.Range("C1").EntireColumn.ColumnWidth = 25.41
$wsAlerts.Range("D1").EntireColumn.ColumnWidth = 18.29
$wsAlerts.Range("E1").EntireColumn.ColumnWidth = 21.14
$wsAlerts.Range("F1").EntireColumn.ColumnWidth = 20.14
$wsAlerts.Range("G1").EntireColumn.ColumnWidth = 15.87

$wsAlerts.Range("C2").Select()

# ---------------------------------------------------------------------------
# New RecordCreditNote sheet, appended after Alerts and made the active tab
# ---------------------------------------------------------------------------
$wsRCN = $wb.Worksheets.Add()
$wsRCN.Name = "RecordCreditNote"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRCN.Move($null, $lastSheet)
$wsRCN = $wb.Worksheets.Item("RecordCreditNote")

$wsRCN.Range("A1").Value = "RECORDCREDITNOTE_UBIN_ITEMNUMBER"
$wsRCN.Range("B1").Value = "RECORDCREDITNOTE_WEARER_NAME"
$wsRCN.Range("C1").Value = "RECORDCREDITNOTE_NOTES"
$wsRCN.Range("D1").Value = "RECORDCREDITNOTE_OBSERVATION"
$wsRCN.Range("E1").Value = "RECORDCREDITNOTE_START_DATE_FORMAT"
$wsRCN.Range("F1").Value = "RECORDCREDITNOTE_END_DATE_FORMAT"

$wsRCN.Range("A2").Value = "Testing"
$wsRCN.Range("B2").Value = "Wearer name from excel"
$wsRCN.Range("C2").Value = "Credit notes from excel sheet"
$wsRCN.Range("D2").Value = "Observation from excel"
$wsRCN.Range("E2").NumberFormat = "DD/MM/YY"
$wsRCN.Range("E2").Value = "16/10/2020"
$wsRCN.Range("F2").NumberFormat = "DD/MM/YY"
$wsRCN.Range("F2").Value = "17/10/2020"

$wsRCN.Range("A1:F1").EntireRow.RowHeight = 15
$wsRCN.Range("A1").EntireColumn.ColumnWidth = 37.42
$wsRCN.Range("B1").EntireColumn.ColumnWidth = 35.13
$wsRCN.Range("C1").EntireColumn.ColumnWidth = 26.29
$wsRCN.Range("D1").EntireColumn.ColumnWidth = 33.57
$wsRCN.Range("E1").EntireColumn.ColumnWidth = 38.29
$wsRCN.Range("F1").EntireColumn.ColumnWidth = 37.92

$wsRCN.Range("C2").Select()
$wsRCN.Select()
